# Applies targeted cell-value updates to the Leve profit-tracking sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) produced by the scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 6000.5
$ws.Range("I20").Value = 6000.5
$ws.Range("K20").Value = 6000.5
$ws.Range("M20").Value = -5770.5
$ws.Range("H35").Value = 6000.5
$ws.Range("I35").Value = 6000.5
$ws.Range("K35").Value = 6000.5
$ws.Range("M35").Value = -5621.5
$ws.Range("H76").Value = 3170
$ws.Range("I76").Value = 3163.158
$ws.Range("J76").Value = 3300
$ws.Range("K76").Value = 3163.158
$ws.Range("L76").Value = 3300
$ws.Range("M76").Value = -2848.158
$ws.Range("N76").Value = -3930
$ws.Range("H79").Value = 3170
$ws.Range("I79").Value = 3163.158
$ws.Range("J79").Value = 3300
$ws.Range("K79").Value = 3163.158
$ws.Range("L79").Value = 3300
$ws.Range("M79").Value = -2071.158
$ws.Range("N79").Value = -5484
$ws.Range("H116").Value = 6162.1
$ws.Range("I116").Value = 1102.5
$ws.Range("J116").Value = 26400.5
$ws.Range("K116").Value = 1102.5
$ws.Range("L116").Value = 26400.5
$ws.Range("M116").Value = 2339.5
$ws.Range("N116").Value = -33284.5
$ws.Range("H118").Value = 666.6667
$ws.Range("I118").Value = 450
$ws.Range("K118").Value = 1350
$ws.Range("M118").Value = 307
$ws.Range("H132").Value = 1726.7778
$ws.Range("I132").Value = 728.65216
$ws.Range("J132").Value = 7466
$ws.Range("K132").Value = 2185.95648
$ws.Range("L132").Value = 22398
$ws.Range("M132").Value = 344.0435200000002
$ws.Range("N132").Value = -27458
$ws.Range("H133").Value = 56280
$ws.Range("J133").Value = 56280
$ws.Range("L133").Value = 56280
$ws.Range("N133").Value = -66400

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 3471
$ws.Range("I31").Value = 3471
$ws.Range("K31").Value = 3471
$ws.Range("M31").Value = -3177
$ws.Range("H35").Value = 50000
$ws.Range("J35").Value = 50000
$ws.Range("L35").Value = 50000
$ws.Range("N35").Value = -50812
$ws.Range("H93").Value = 30789
$ws.Range("J93").Value = 30789
$ws.Range("L93").Value = 30789
$ws.Range("N93").Value = -35781
$ws.Range("H132").Value = 2418.7
$ws.Range("I132").Value = 2164.2058
$ws.Range("J132").Value = 2959.5
$ws.Range("K132").Value = 6492.617400000001
$ws.Range("L132").Value = 8878.5
$ws.Range("M132").Value = -3962.617400000001
$ws.Range("N132").Value = -13938.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 21900
$ws.Range("I52").Value = 5700
$ws.Range("J52").Value = 30000
$ws.Range("K52").Value = 5700
$ws.Range("L52").Value = 30000
$ws.Range("M52").Value = -5437
$ws.Range("N52").Value = -30526
$ws.Range("H94").Value = 1534.9166
$ws.Range("J94").Value = 1224
$ws.Range("L94").Value = 1224
$ws.Range("N94").Value = -2126
$ws.Range("H99").Value = 1633.3334
$ws.Range("I99").Value = 1383.3334
$ws.Range("K99").Value = 1383.3334
$ws.Range("M99").Value = 114.6666
$ws.Range("H106").Value = 60671
$ws.Range("J106").Value = 60671
$ws.Range("L106").Value = 60671
$ws.Range("N106").Value = -63195
$ws.Range("H109").Value = 27998
$ws.Range("J109").Value = 27998
$ws.Range("L109").Value = 27998
$ws.Range("N109").Value = -30772
$ws.Range("H118").Value = 35000
$ws.Range("J118").Value = 35000
$ws.Range("L118").Value = 35000
$ws.Range("N118").Value = -38314
$ws.Range("H121").Value = 21900
$ws.Range("I121").Value = 5700
$ws.Range("J121").Value = 30000
$ws.Range("K121").Value = 5700
$ws.Range("L121").Value = 30000
$ws.Range("M121").Value = -3953
$ws.Range("N121").Value = -33494
$ws.Range("H126").Value = 35000
$ws.Range("J126").Value = 35000
$ws.Range("L126").Value = 35000
$ws.Range("N126").Value = -44880
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 975.7857
$ws.Range("I16").Value = 938.4167
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 938.4167
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -651.4167
$ws.Range("N16").Value = -1774
$ws.Range("H58").Value = 1017.6087
$ws.Range("I58").Value = 880.25
$ws.Range("J58").Value = 1933.3334
$ws.Range("K58").Value = 880.25
$ws.Range("L58").Value = 1933.3334
$ws.Range("M58").Value = -677.25
$ws.Range("N58").Value = -2339.3334
$ws.Range("H107").Value = 2059.7896
$ws.Range("I107").Value = 421.5
$ws.Range("J107").Value = 3251.2727
$ws.Range("K107").Value = 421.5
$ws.Range("L107").Value = 3251.2727
$ws.Range("M107").Value = 1498.5
$ws.Range("N107").Value = -7091.2727
$ws.Range("H113").Value = 975.7857
$ws.Range("I113").Value = 938.4167
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 938.4167
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1231.5833
$ws.Range("N113").Value = -5540
$ws.Range("H130").Value = 20000
$ws.Range("J130").Value = 20000
$ws.Range("L130").Value = 20000
$ws.Range("N130").Value = -30040
$ws.Range("H132").Value = 2088.75
$ws.Range("I132").Value = 1864.8846
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 5594.6538
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -3064.6538
$ws.Range("N132").Value = -20057
$ws.Range("H134").Value = 31430044
$ws.Range("I134").Value = 3031757.8
$ws.Range("J134").Value = 500001760
$ws.Range("K134").Value = 9095273.399999999
$ws.Range("L134").Value = 1500005280
$ws.Range("M134").Value = -9092738.399999999
$ws.Range("N134").Value = -1500010350
$ws.Range("H136").Value = 1017.6087
$ws.Range("I136").Value = 880.25
$ws.Range("J136").Value = 1933.3334
$ws.Range("K136").Value = 2640.75
$ws.Range("L136").Value = 5800.0002
$ws.Range("M136").Value = -90.75
$ws.Range("N136").Value = -10900.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1254.5454
$ws.Range("I98").Value = 1488
$ws.Range("K98").Value = 4464
$ws.Range("M98").Value = -2966

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 27333.334
$ws.Range("J23").Value = 27333.334
$ws.Range("L23").Value = 27333.334
$ws.Range("N23").Value = -27779.334
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H113").Value = 1345.3334
$ws.Range("I113").Value = 1032.3125
$ws.Range("J113").Value = 2347
$ws.Range("K113").Value = 1032.3125
$ws.Range("L113").Value = 2347
$ws.Range("M113").Value = 1137.6875
$ws.Range("N113").Value = -6687
$ws.Range("H128").Value = 20000
$ws.Range("J128").Value = 20000
$ws.Range("L128").Value = 20000
$ws.Range("N128").Value = -29960
$ws.Range("H130").Value = 216000
$ws.Range("J130").Value = 216000
$ws.Range("L130").Value = 216000
$ws.Range("N130").Value = -226040
$ws.Range("H132").Value = 2715.2195
$ws.Range("I132").Value = 2716.9
$ws.Range("J132").Value = 2710.6365
$ws.Range("K132").Value = 8150.700000000001
$ws.Range("L132").Value = 8131.9095
$ws.Range("M132").Value = -5620.700000000001
$ws.Range("N132").Value = -13191.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 31892.166
$ws.Range("J112").Value = 31892.166
$ws.Range("L112").Value = 31892.166
$ws.Range("N112").Value = -34846.166
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H121").Value = 28842
$ws.Range("J121").Value = 28842
$ws.Range("L121").Value = 28842
$ws.Range("N121").Value = -32336
$ws.Range("H122").Value = 5004.727
$ws.Range("I122").Value = 5673.6
$ws.Range("K122").Value = 17020.8
$ws.Range("M122").Value = -14570.8
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H127").Value = 34992.23
$ws.Range("J127").Value = 34992.23
$ws.Range("L127").Value = 34992.23
$ws.Range("N127").Value = -44912.23
$ws.Range("H132").Value = 3890.3667
$ws.Range("I132").Value = 3856.0952
$ws.Range("J132").Value = 3970.3333
$ws.Range("K132").Value = 11568.2856
$ws.Range("L132").Value = 11910.9999
$ws.Range("M132").Value = -9038.285600000001
$ws.Range("N132").Value = -16970.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 11366.444
$ws.Range("J69").Value = 11366.444
$ws.Range("L69").Value = 11366.444
$ws.Range("N69").Value = -12864.444
$ws.Range("H72").Value = 11366.444
$ws.Range("J72").Value = 11366.444
$ws.Range("L72").Value = 34099.33199999999
$ws.Range("N72").Value = -41587.33199999999
$ws.Range("H107").Value = 5038.25
$ws.Range("I107").Value = 8100.643
$ws.Range("K107").Value = 24301.929
$ws.Range("M107").Value = -22381.929
$ws.Range("H113").Value = 486.86667
$ws.Range("I113").Value = 380.33334
$ws.Range("J113").Value = 646.6667
$ws.Range("K113").Value = 1141.00002
$ws.Range("L113").Value = 1940.0001
$ws.Range("M113").Value = 1028.99998
$ws.Range("N113").Value = -6280.0001
